$d = $word.ActiveDocument

# Locate the first bullet paragraph of the "KEY ACHIEVEMENTS AND IMPACT"
# section by its distinctive trailing text (unique within the document,
# unlike the shorter variant of this sentence used earlier in the
# "Partner - Siege Analytics" experience bullets).
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*serving 12,847 analysts across 89 organizations*") {
        $p1 = $p
        break
    }
}

# Walk forward to grab the rest of the six bullets making up this section.
$p2 = $p1.Next()
$p3 = $p2.Next()
$p4 = $p3.Next()
$p5 = $p4.Next()
$p6 = $p5.Next()

# Rewrite the first three bullets as impact-focused accomplishment
# statements.
$p1.Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$p2.Range.Text = "• `$4.7M savings enabled nonprofit access"
$p3.Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"

# Collapse the remaining three bullets into a single accomplishment
# statement: rewrite the fourth in place, then delete the fifth and sixth
# paragraphs entirely. Delete from the end backwards so that removing a
# later paragraph's range doesn't shift/invalidate the still-pending
# paragraph reference ahead of it.
$p4.Range.Text = "• Real-time collaboration at national scale"
$p6.Range.Delete()
$p5.Range.Delete()
